$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M5").Value = 70
$ws.Range("I5").Value = 45
$ws.Range("K5").Value = 45
$ws.Range("H5").Value = 45
$ws.Range("M69").Value = -24140
$ws.Range("H69").Value = 9816.929
$ws.Range("K69").Value = 25014
$ws.Range("I69").Value = 8338
$ws.Range("I72").Value = 8338
$ws.Range("K72").Value = 75042
$ws.Range("M72").Value = -70674
$ws.Range("H72").Value = 9816.929
$ws.Range("L105").Value = 68992
$ws.Range("H105").Value = 68992
$ws.Range("J105").Value = 68992
$ws.Range("N105").Value = -75980
$ws.Range("L117").Value = 68948.25
$ws.Range("J117").Value = 68948.25
$ws.Range("N117").Value = -78126.25
$ws.Range("H117").Value = 68948.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K26").Value = 2966.3333
$ws.Range("I26").Value = 2966.3333
$ws.Range("M26").Value = -2636.3333
$ws.Range("H26").Value = 2966.3333
$ws.Range("K63").Value = 2124.5
$ws.Range("H63").Value = 3797.4546
$ws.Range("I63").Value = 2124.5
$ws.Range("M63").Value = -1438.5
$ws.Range("I66").Value = 2124.5
$ws.Range("H66").Value = 3797.4546
$ws.Range("K66").Value = 10622.5
$ws.Range("M66").Value = -7190.5
$ws.Range("M88").Value = -558.1429000000001
$ws.Range("H88").Value = 1077.5834
$ws.Range("I88").Value = 964.1429000000001
$ws.Range("N88").Value = -2048.4
$ws.Range("K88").Value = 964.1429000000001
$ws.Range("J88").Value = 1236.4
$ws.Range("L88").Value = 1236.4
$ws.Range("L91").Value = 1236.4
$ws.Range("N91").Value = -4044.4
$ws.Range("M91").Value = 439.8570999999999
$ws.Range("J91").Value = 1236.4
$ws.Range("I91").Value = 964.1429000000001
$ws.Range("H91").Value = 1077.5834
$ws.Range("K91").Value = 964.1429000000001
$ws.Range("N124").Value = -79820
$ws.Range("H124").Value = 70000
$ws.Range("J124").Value = 70000
$ws.Range("L124").Value = 70000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N40").Value = -59419.668
$ws.Range("J40").Value = 58889.668
$ws.Range("H40").Value = 58889.668
$ws.Range("L40").Value = 58889.668
$ws.Range("K96").Value = 28333
$ws.Range("N96").Value = -99446
$ws.Range("H96").Value = 54581.4
$ws.Range("I96").Value = 28333
$ws.Range("L96").Value = 93954
$ws.Range("J96").Value = 93954
$ws.Range("M96").Value = -25587
$ws.Range("L125").Value = 130994.5
$ws.Range("H125").Value = 130994.5
$ws.Range("N125").Value = -140834.5
$ws.Range("J125").Value = 130994.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L18").Value = 100990
$ws.Range("J18").Value = 100990
$ws.Range("H18").Value = 100990
$ws.Range("N18").Value = -101450
$ws.Range("K58").Value = 2478.8572
$ws.Range("H58").Value = 2696.309
$ws.Range("I58").Value = 2478.8572
$ws.Range("M58").Value = -2275.8572
$ws.Range("N87").Value = -84943.39999999999
$ws.Range("L87").Value = 82571.39999999999
$ws.Range("H87").Value = 82571.39999999999
$ws.Range("J87").Value = 82571.39999999999
$ws.Range("N90").Value = -259570.2
$ws.Range("H90").Value = 82571.39999999999
$ws.Range("L90").Value = 247714.2
$ws.Range("J90").Value = 82571.39999999999
$ws.Range("N124").Value = -89909
$ws.Range("H124").Value = 84999
$ws.Range("J124").Value = 84999
$ws.Range("L124").Value = 84999
$ws.Range("K132").Value = 9477.332999999999
$ws.Range("I132").Value = 3159.111
$ws.Range("M132").Value = -6947.332999999999
$ws.Range("H132").Value = 4067.2273
$ws.Range("I134").Value = 1686.9286
$ws.Range("L134").Value = 7398
$ws.Range("J134").Value = 2466
$ws.Range("K134").Value = 5060.7858
$ws.Range("N134").Value = -12468
$ws.Range("M134").Value = -2525.7858
$ws.Range("H134").Value = 1920.65
$ws.Range("H136").Value = 2696.309
$ws.Range("M136").Value = -4886.571599999999
$ws.Range("K136").Value = 7436.571599999999
$ws.Range("I136").Value = 2478.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J12").Value = 769331.9
$ws.Range("K12").Value = 72
$ws.Range("N12").Value = -2308341.7
$ws.Range("M12").Value = 101
$ws.Range("I12").Value = 24
$ws.Range("H12").Value = 666757.4399999999
$ws.Range("L12").Value = 2307995.7
$ws.Range("N68").Value = -4374.8
$ws.Range("H68").Value = 898.3333
$ws.Range("J68").Value = 917.6
$ws.Range("I68").Value = 802
$ws.Range("K68").Value = 2406
$ws.Range("L68").Value = 2752.8
$ws.Range("M68").Value = -1595
$ws.Range("L71").Value = 8258.4
$ws.Range("I71").Value = 802
$ws.Range("M71").Value = -3162
$ws.Range("K71").Value = 7218
$ws.Range("N71").Value = -16370.4
$ws.Range("H71").Value = 898.3333
$ws.Range("J71").Value = 917.6
$ws.Range("H74").Value = 15977.5
$ws.Range("J74").Value = 20636.666
$ws.Range("N74").Value = -64031.99800000001
$ws.Range("L74").Value = 61909.99800000001
$ws.Range("L77").Value = 185729.994
$ws.Range("J77").Value = 20636.666
$ws.Range("H77").Value = 15977.5
$ws.Range("N77").Value = -196337.994
$ws.Range("N114").Value = -14616.0001
$ws.Range("L114").Value = 8108.000100000001
$ws.Range("J114").Value = 2702.6667
$ws.Range("H114").Value = 1965.7142
$ws.Range("N129").Value = -17798.7145
$ws.Range("L129").Value = 7798.7145
$ws.Range("H129").Value = 1765.6364
$ws.Range("J129").Value = 2599.5715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("J119").Value = 69800.39999999999
$ws.Range("L119").Value = 69800.39999999999
$ws.Range("N119").Value = -79476.39999999999
$ws.Range("H119").Value = 69800.39999999999
$ws.Range("N124").Value = -87728.5
$ws.Range("H124").Value = 77908.5
$ws.Range("J124").Value = 77908.5
$ws.Range("L124").Value = 77908.5
$ws.Range("K132").Value = 8726.143199999999
$ws.Range("I132").Value = 2908.7144
$ws.Range("M132").Value = -6196.143199999999
$ws.Range("H132").Value = 3335.1853
$ws.Range("N80").ClearContents()
$ws.Range("M80").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N40").Value = -3247
$ws.Range("I40").Value = 2474.75
$ws.Range("J40").Value = 2975
$ws.Range("M40").Value = -2338.75
$ws.Range("H40").Value = 2752.6667
$ws.Range("K40").Value = 2474.75
$ws.Range("L40").Value = 2975
$ws.Range("N92").Value = -113991
$ws.Range("H92").Value = 108999
$ws.Range("J92").Value = 108999
$ws.Range("L92").Value = 108999
$ws.Range("N122").Value = -47802.625
$ws.Range("H122").Value = 12778.695
$ws.Range("L122").Value = 42902.625
$ws.Range("M122").Value = -33450.598
$ws.Range("K122").Value = 35900.598
$ws.Range("I122").Value = 11966.866
$ws.Range("J122").Value = 14300.875
$ws.Range("H136").Value = 5255.5186
$ws.Range("N136").Value = -30294
$ws.Range("M136").Value = -7671
$ws.Range("L136").Value = 25194
$ws.Range("K136").Value = 10221
$ws.Range("I136").Value = 3407
$ws.Range("J136").Value = 8398

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L57").Value = 128659.664
$ws.Range("N57").Value = -130167.664
$ws.Range("J57").Value = 128659.664
$ws.Range("H57").Value = 128659.664
$ws.Range("M81").Value = -2366.8572
$ws.Range("J81").Value = 5178.9
$ws.Range("N81").Value = -12479.8
$ws.Range("K81").Value = 3427.8572
$ws.Range("I81").Value = 1713.9286
$ws.Range("H81").Value = 3157.6667
$ws.Range("L81").Value = 10357.8
$ws.Range("I84").Value = 1713.9286
$ws.Range("L84").Value = 51789
$ws.Range("J84").Value = 5178.9
$ws.Range("N84").Value = -62397
$ws.Range("K84").Value = 17139.286
$ws.Range("H84").Value = 3157.6667
$ws.Range("M84").Value = -11835.286
$ws.Range("H126").Value = 2204.1667
$ws.Range("M126").Value = 13.75
$ws.Range("K126").Value = 2456.25
$ws.Range("I126").Value = 818.75
$ws.Range("H136").Value = 1464.7037
$ws.Range("M136").Value = -1417.2498
$ws.Range("K136").Value = 3967.2498
$ws.Range("I136").Value = 1322.4166
